# Applies the recorded "output generated at 456a3b4" update to
# 上海-漫展信息.xlsx across its four worksheets:
#   1. 展览      (Exhibitions)
#   2. 演出      (Performances)
#   3. 本地生活   (Local life)
#   4. 全部类型   (All types)

$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet 1: 展览 (Exhibitions) - "want to go" count (column F) bumps
# =====================================================================
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 679
$ws1.Range("F4").Value = 3144
$ws1.Range("F8").Value = 317
$ws1.Range("F9").Value = 7370
$ws1.Range("F13").Value = 438
$ws1.Range("F15").Value = 1943
$ws1.Range("F16").Value = 1796
$ws1.Range("F18").Value = 26
$ws1.Range("F19").Value = 74
$ws1.Range("F20").Value = 1849
$ws1.Range("F21").Value = 1381
$ws1.Range("F22").Value = 1257
$ws1.Range("F23").Value = 653
$ws1.Range("F25").Value = 1143
$ws1.Range("F27").Value = 545
$ws1.Range("F28").Value = 133
$ws1.Range("F29").Value = 75
$ws1.Range("F30").Value = 4709
$ws1.Range("F31").Value = 2277
$ws1.Range("F32").Value = 3910
$ws1.Range("F33").Value = 2194
$ws1.Range("F34").Value = 161
$ws1.Range("F35").Value = 226
$ws1.Range("F36").Value = 1153
$ws1.Range("F38").Value = 50
$ws1.Range("F40").Value = 383
$ws1.Range("F42").Value = 158
$ws1.Range("F43").Value = 525
$ws1.Range("F44").Value = 258
$ws1.Range("F46").Value = 793
$ws1.Range("F48").Value = 6
$ws1.Range("F49").Value = 158

# =====================================================================
# Sheet 2: 演出 (Performances)
# =====================================================================
$ws2 = $wb.Worksheets.Item("演出")

# -- simple "want to go" bumps --
$ws2.Range("F14").Value = 95
$ws2.Range("F17").Value = 551
$ws2.Range("F24").Value = 91

# -- row 9: lowest price switched from a numeric value to "sold out"/
#    "not sellable" text --
$ws2.Range("G9").Value = "不可售"

# -- a new event ("生命之舞" - Paul Taylor Dance Company) was inserted
#    as the new row 33, pushing the previous rows 33-35 down to 34-36 --
$ws2.Rows.Item(33).Insert()

# Copy the row-label cell's formatting (bold / border / centered) from
# the row above so the new A33 matches the rest of column A.
$ws2.Range("A32").Copy()
$ws2.Range("A33").PasteSpecial(-4122)

$ws2.Range("A33").Value = 32
$ws2.Range("B33").NumberFormat = "@"
$ws2.Range("B33").Value = "2024-12-05"
$ws2.Range("C33").Value = [char]0x4E0A + [char]0x6D77 + [char]0xB7 + '"' + [char]0x751F + [char]0x547D + [char]0x4E4B + [char]0x821E + '"' + [char]0x4FDD + [char]0x7F57 + [char]0x6CF0 + [char]0x52D2 + [char]0x73B0 + [char]0x4EE3 + [char]0x821E + [char]0x56E2 + [char]0x4E16 + [char]0x7EAA + [char]0x4E4B + [char]0x591C
$ws2.Range("D33").Value = [char]0x4E1C + [char]0x5927 + [char]0x540D + [char]0x8DEF + '889' + [char]0x53F7 + ' ' + [char]0x5317 + [char]0x5916 + [char]0x6EE9 + [char]0x53CB + [char]0x90A6 + [char]0x5927 + [char]0x5267 + [char]0x9662
$ws2.Range("E33").NumberFormat = "@"
$ws2.Range("E33").Value = "2024.12.05 19:30-12.06 22:00"
$ws2.Range("F33").Value = 0
$ws2.Range("G33").Value = 180
$ws2.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=91132"
$ws2.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202408/6ozgDFcQ1724315518809.jpeg"

# -- the (now shifted) rows 35 and 36 also got their "want to go" counts
#    bumped by one --
$ws2.Range("F35").Value = 45
$ws2.Range("F36").Value = 18

# =====================================================================
# Sheet 3: 本地生活 (Local life) - "want to go" count bumps
# =====================================================================
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 564
$ws3.Range("F6").Value = 1853
$ws3.Range("F7").Value = 1885
$ws3.Range("F8").Value = 2903
$ws3.Range("F9").Value = 1145
$ws3.Range("F10").Value = 1136
$ws3.Range("F12").Value = 438
$ws3.Range("F13").Value = 1844
$ws3.Range("F14").Value = 8217
$ws3.Range("F15").Value = 360

# =====================================================================
# Sheet 4: 全部类型 (All types) - "want to go" count bumps
# =====================================================================
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 679
$ws4.Range("F4").Value = 3144
$ws4.Range("F6").Value = 1853
$ws4.Range("F7").Value = 317
$ws4.Range("F8").Value = 2903
$ws4.Range("F9").Value = 7370
$ws4.Range("F10").Value = 1145
$ws4.Range("F11").Value = 1136
$ws4.Range("F13").Value = 438
$ws4.Range("F15").Value = 438
$ws4.Range("F18").Value = 26
$ws4.Range("F19").Value = 74
$ws4.Range("F20").Value = 1849
$ws4.Range("F21").Value = 1381
$ws4.Range("F22").Value = 1257
$ws4.Range("F23").Value = 653
$ws4.Range("F25").Value = 1143
$ws4.Range("F26").Value = 95
$ws4.Range("F29").Value = 551
$ws4.Range("F30").Value = 545
$ws4.Range("F32").Value = 133
$ws4.Range("F33").Value = 75
$ws4.Range("F34").Value = 4709
$ws4.Range("F35").Value = 2277
$ws4.Range("F36").Value = 3910
$ws4.Range("F37").Value = 2194
$ws4.Range("F38").Value = 162
$ws4.Range("F39").Value = 226
$ws4.Range("F40").Value = 1153
$ws4.Range("F43").Value = 383
$ws4.Range("F44").Value = 158
$ws4.Range("F45").Value = 91
$ws4.Range("F46").Value = 525
$ws4.Range("F47").Value = 258

# -- row 16: the "1PLUS1" pop-up market event was replaced by a
#    "VWonderland" vtuber meet-and-greet event --
$ws4.Range("C16").Value = [char]0x4E0A + [char]0x6D77 + [char]0xB7 + 'VWonderland' + [char]0x865A + [char]0x62DF + [char]0x4E3B + [char]0x64AD + [char]0x7EBF + [char]0x4E0B + [char]0x89C1 + [char]0x9762 + [char]0x4F1A
$ws4.Range("D16").Value = [char]0x7FD4 + [char]0x6BB7 + [char]0x8DEF + '1099' + [char]0x53F7 + ' ' + [char]0x5408 + [char]0x751F + [char]0x6C47
$ws4.Range("E16").NumberFormat = "@"
$ws4.Range("E16").Value = "2024.08.24 10:00-08.24 21:00"
$ws4.Range("F16").Value = 86
$ws4.Range("G16").Value = 60
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=90693"
$ws4.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202408/FZ9CsGO81723560782092.png"

Write-Output "edit complete"
